$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 311
$ws.Range("J17").Value = 264.90244
$ws.Range("L17").Value = 794.70732
$ws.Range("N17").Value = -1130.70732

$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("N41").Value = 0

$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("N51").Value = 0

$ws.Range("H62").Value = 1721.8334
$ws.Range("I62").Value = 1066.2
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 1066.2
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -442.2
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 1721.8334
$ws.Range("I65").Value = 1066.2
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 5331
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -2211
$ws.Range("N65").Value = -31240

$ws.Range("H92").Value = 5291217
$ws.Range("I92").Value = 223.6842
$ws.Range("K92").Value = 223.6842
$ws.Range("M92").Value = 1024.3158

$ws.Range("H112").Value = 1327
$ws.Range("I112").Value = 1333.3334
$ws.Range("J112").Value = 1326.174
$ws.Range("K112").Value = 4000.0002
$ws.Range("L112").Value = 3978.522
$ws.Range("M112").Value = -2892.0002
$ws.Range("N112").Value = -6194.522

$ws.Range("H116").Value = 7958.2354
$ws.Range("I116").Value = 5260
$ws.Range("K116").Value = 5260
$ws.Range("M116").Value = -1818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 28638.54
$ws.Range("I2").Value = 38480.223
$ws.Range("K2").Value = 38480.223
$ws.Range("M2").Value = -38367.223

$ws.Range("H32").Value = 4042.257
$ws.Range("I32").Value = 4300.361
$ws.Range("K32").Value = 4300.361
$ws.Range("M32").Value = -4013.361

$ws.Range("H45").Value = 8398.799999999999
$ws.Range("I45").Value = 4498
$ws.Range("K45").Value = 4498
$ws.Range("M45").Value = -4121

$ws.Range("H63").Value = 6969.391
$ws.Range("I63").Value = 5164.2144
$ws.Range("J63").Value = 9777.444
$ws.Range("K63").Value = 5164.2144
$ws.Range("L63").Value = 9777.444
$ws.Range("M63").Value = -4478.2144
$ws.Range("N63").Value = -11149.444

$ws.Range("H66").Value = 6969.391
$ws.Range("I66").Value = 5164.2144
$ws.Range("J66").Value = 9777.444
$ws.Range("K66").Value = 25821.072
$ws.Range("L66").Value = 48887.22
$ws.Range("M66").Value = -22389.072
$ws.Range("N66").Value = -55751.22

$ws.Range("H74").Value = 4966.067
$ws.Range("I74").Value = 4457.5835
$ws.Range("J74").Value = 7000
$ws.Range("K74").Value = 4457.5835
$ws.Range("L74").Value = 7000
$ws.Range("M74").Value = -3583.5835
$ws.Range("N74").Value = -8748

$ws.Range("H77").Value = 4966.067
$ws.Range("I77").Value = 4457.5835
$ws.Range("J77").Value = 7000
$ws.Range("K77").Value = 22287.9175
$ws.Range("L77").Value = 35000
$ws.Range("M77").Value = -17919.9175
$ws.Range("N77").Value = -43736

$ws.Range("H116").Value = 28638.54
$ws.Range("I116").Value = 38480.223
$ws.Range("K116").Value = 38480.223
$ws.Range("M116").Value = -36186.223

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 28638.54
$ws.Range("I3").Value = 38480.223
$ws.Range("K3").Value = 38480.223
$ws.Range("M3").Value = -38366.223

$ws.Range("H20").Value = 4752.8335
$ws.Range("I20").Value = 3582.3333
$ws.Range("J20").Value = 5923.3335
$ws.Range("K20").Value = 3582.3333
$ws.Range("L20").Value = 5923.3335
$ws.Range("M20").Value = -3335.3333
$ws.Range("N20").Value = -6417.3335

$ws.Range("H94").Value = 3006538.8
$ws.Range("I94").Value = 2818.4138
$ws.Range("J94").Value = 13895025
$ws.Range("K94").Value = 2818.4138
$ws.Range("L94").Value = 13895025
$ws.Range("M94").Value = -2367.4138
$ws.Range("N94").Value = -13895927

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 7227.1177
$ws.Range("I107").Value = 712.1667
$ws.Range("J107").Value = 22863
$ws.Range("K107").Value = 712.1667
$ws.Range("L107").Value = 22863
$ws.Range("M107").Value = 1207.8333
$ws.Range("N107").Value = -26703

$ws.Range("H129").Value = 50000
$ws.Range("J129").Value = 50000
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000

$ws.Range("H132").Value = 2550.5
$ws.Range("I132").Value = 2276.4614
$ws.Range("K132").Value = 6829.3842
$ws.Range("M132").Value = -4299.3842

$ws.Range("H134").Value = 7964
$ws.Range("I134").Value = 7964
$ws.Range("K134").Value = 23892
$ws.Range("M134").Value = -21357

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2021.5
$ws.Range("J5").Value = 1121
$ws.Range("L5").Value = 3363
$ws.Range("N5").Value = -3587

$ws.Range("H132").Value = 1939.9
$ws.Range("I132").Value = 1466.6666
$ws.Range("K132").Value = 13199.9994
$ws.Range("M132").Value = -10669.9994

$ws.Range("H135").Value = 2021.5
$ws.Range("J135").Value = 1121
$ws.Range("L135").Value = 10089
$ws.Range("N135").Value = -15159

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 21969.666
$ws.Range("I33").Value = 20000
$ws.Range("J33").Value = 23939.334
$ws.Range("K33").Value = 20000
$ws.Range("L33").Value = 23939.334
$ws.Range("M33").Value = -19748
$ws.Range("N33").Value = -24443.334

$ws.Range("H38").Value = 21493.5
$ws.Range("J38").Value = 21493.5
$ws.Range("L38").Value = 21493.5
$ws.Range("N38").Value = -22419.5

$ws.Range("H40").Value = 19285.666
$ws.Range("J40").Value = 20199.5
$ws.Range("L40").Value = 20199.5
$ws.Range("N40").Value = -20501.5

$ws.Range("H44").Value = 12886.143
$ws.Range("I44").Value = 10551
$ws.Range("K44").Value = 10551
$ws.Range("M44").Value = -9955

$ws.Range("H47").Value = 25197.75
$ws.Range("I47").Value = 23000
$ws.Range("J47").Value = 25930.334
$ws.Range("K47").Value = 23000
$ws.Range("L47").Value = 25930.334
$ws.Range("M47").Value = -22432
$ws.Range("N47").Value = -27066.334

$ws.Range("H122").Value = 3860.6875
$ws.Range("I122").Value = 3314.6365
$ws.Range("J122").Value = 5062
$ws.Range("K122").Value = 9943.9095
$ws.Range("L122").Value = 15186
$ws.Range("M122").Value = -7493.9095
$ws.Range("N122").Value = -20086

$ws.Range("H129").Value = 69780
$ws.Range("J129").Value = 69780
$ws.Range("L129").Value = 69780
$ws.Range("N129").Value = -79780

$ws.Range("H132").Value = 8363.125
$ws.Range("I132").Value = 7186.269
$ws.Range("J132").Value = 13462.833
$ws.Range("K132").Value = 21558.807
$ws.Range("L132").Value = 40388.499
$ws.Range("M132").Value = -19028.807
$ws.Range("N132").Value = -45448.499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 28374.5
$ws.Range("J110").Value = 28374.5
$ws.Range("L110").Value = 28374.5
$ws.Range("N110").Value = -36554.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5432
$ws.Range("I132").Value = 3562.842
$ws.Range("K132").Value = 10688.526
$ws.Range("M132").Value = -8158.526
